$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.771.69"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "2.293.25"
$ws.Range("E3").Value = "  -1.13%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "102.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "270.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("E7").Value = "  -0.43%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.608"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.07"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0931"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.07%  "
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.856"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.27%  "
$ws.Range("D16").Value = "2.292.91"
$ws.Range("E16").Value = "  -0.82%  "
$ws.Range("D17").Value = "43.763.89"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("E19").Value = "  -2.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.33%  "
$ws.Range("E21").Value = "  +10.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "233.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E24").Value = "  -1.99%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "40.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.93%  "
$ws.Range("E28").Value = "  -1.55%  "
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "177.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "21.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0904"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.43%  "
$ws.Range("E34").Value = "  +9.75%  "
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0359"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.56"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.65%  "
$ws.Range("E39").Value = "  -2.39%  "
$ws.Range("E40").Value = "  -2.36%  "
$ws.Range("E41").Value = "  -2.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.32%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E45").Value = "  -4.61%  "
$ws.Range("E46").Value = "  -0.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "99.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.40%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +11.70%  "
$ws.Range("B50").Value = "WOONetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.442"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.54%  "
$ws.Range("D51").Value = "2.514.34"
$ws.Range("E51").Value = "  -1.17%  "
